# Update Bakery (cake) prices on Sheet1 and refresh the current selection
# to match where the author was working when the file was re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D24").Value = 650
$ws.Range("D25").Value = 550
$ws.Range("D26").Value = 600
$ws.Range("D27").Value = 400
$ws.Range("D28").Value = 550
$ws.Range("D29").Value = 500

# Move the visible selection to D25, matching the saved workbook state.
[void]$ws.Range("D25").Select()
